$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date/time number format from the cell above (G5) onto the new
# G6 cell so it reuses the existing style (numFmtId 22, "m/d/yyyy h:mm")
# instead of Excel creating a brand new custom number format.
$ws.Range("G5").Copy() | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Append the new trade row (row 6) with the same columns as the rows above it.
$ws.Range("A6").Value = 10156.799999999999
$ws.Range("B6").Value = 10107.27
$ws.Range("C6").Value = 307.20999999999998
$ws.Range("D6").Value = 308.70999999999998
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = 0.49
$ws.Range("G6").Value = 42609.503946759258
$ws.Range("H6").Value = $true
